# Refresh the "ランサーズ" (Lancers) listing sheet:
#  - new scrape timestamp for every remaining row
#  - rows 2-3 keep their listings, rows 4-12 get new listing data
#    (old row 9's listing moves up to row 4, several rows are brand new,
#    and a few more old rows shift up into rows 7/11/12)
#  - old rows 13-20 are dropped entirely (sheet shrinks from H20 to H12)
#  - column D narrows from 32 to 30 characters wide
#  - the F-column hyperlinks are rebuilt so they keep pointing at the
#    same URL that is now displayed in each cell

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Drop the old trailing rows (13-20) first so everything below stays
#    simple arithmetic on row numbers for the remaining edits.
# ---------------------------------------------------------------------
$ws.Rows("13:20").Delete()

# ---------------------------------------------------------------------
# 2) Clear out the stale hyperlink objects. (Row-delete above leaves the
#    worksheet's Hyperlinks collection referencing the removed rows, and
#    we are about to change several F-column URLs anyway, so just drop
#    all of them and rebuild only the ones we need.)
# ---------------------------------------------------------------------
$ws.Hyperlinks.Delete()

# ---------------------------------------------------------------------
# 3) New scrape timestamp, identical for every data row.
# ---------------------------------------------------------------------
$timestamp = "2025-12-25 06:30:05"
for ($r = 2; $r -le 12; $r++) {
    $ws.Cells.Item($r, 1).Value = $timestamp
}

# ---------------------------------------------------------------------
# 4) Row-by-row content (B title, C category, D price, E deadline,
#    F URL, G score, H skills). The F-column URL is kept in its own
#    PowerShell variable so it can be reused verbatim for the matching
#    Hyperlinks.Add() call below (reading .Value back off the cell for
#    that call is unreliable in this runtime).
# ---------------------------------------------------------------------

# Row 2 - unchanged listing, only the timestamp (above) changed.
$ws.Cells.Item(2, 2).Value = "製造業向け図面自動生成システムの開発・ツール化を支援してくださるエンジニア募集(AI/バックエンド)"
$ws.Cells.Item(2, 3).Value = "システム開発"
$ws.Cells.Item(2, 4).Value = "300,000 円 ~ 500,000 円 / 固定"
$ws.Cells.Item(2, 5).Value = "期限情報なし"
$url2 = "https://www.lancers.jp/work/detail/5460562"
$ws.Cells.Item(2, 6).Value = $url2
$ws.Cells.Item(2, 7).Value = 435
$ws.Cells.Item(2, 8).Value = "🔥AI,Ai ◆ツール,開発"

# Row 3 - unchanged listing, only the timestamp (above) changed.
$ws.Cells.Item(3, 2).Value = "既存の情報検索システム(PHP)にAI文書作成システム(既存システムへの機能追加)の開発者募集します"
$ws.Cells.Item(3, 3).Value = "システム開発"
$ws.Cells.Item(3, 4).Value = "100,000 円 ~ 200,000 円 / 固定"
$ws.Cells.Item(3, 5).Value = "期限情報なし"
$url3 = "https://www.lancers.jp/work/detail/5460357"
$ws.Cells.Item(3, 6).Value = $url3
$ws.Cells.Item(3, 7).Value = 388
$ws.Cells.Item(3, 8).Value = "🔥AI,Ai ◆開発 ○PHP"

# Row 4 - new listing (was previously at row 9).
$ws.Cells.Item(4, 2).Value = "施設管理・現場業務向け チェックリスト業務の自動化・報告書作成システム開発エンジニア募集"
$ws.Cells.Item(4, 3).Value = "システム開発"
$ws.Cells.Item(4, 4).Value = "300,000 円 ~ 500,000 円 / 固定"
$ws.Cells.Item(4, 5).Value = "期限情報なし"
$url4 = "https://www.lancers.jp/work/detail/5460563"
$ws.Cells.Item(4, 6).Value = $url4
$ws.Cells.Item(4, 7).Value = 220
$ws.Cells.Item(4, 8).Value = "◆開発,システム開発 ◇管理"

# Row 5 - brand new listing.
$ws.Cells.Item(5, 2).Value = "【急募】魅力的なECサイトのWebシステム開発依頼"
$ws.Cells.Item(5, 3).Value = "システム開発"
$ws.Cells.Item(5, 4).Value = "200,000 円 ~ 300,000 円 / 固定"
$ws.Cells.Item(5, 5).Value = "期限情報なし"
$url5 = "https://www.lancers.jp/work/detail/5460750"
$ws.Cells.Item(5, 6).Value = $url5
$ws.Cells.Item(5, 7).Value = 153
$ws.Cells.Item(5, 8).Value = "◆開発,システム開発 ◇サイト"

# Row 6 - brand new listing.
$ws.Cells.Item(6, 2).Value = "【急募】魅力的なWebシステム開発の提案をお待ちしています!"
$ws.Cells.Item(6, 3).Value = "システム開発"
$ws.Cells.Item(6, 4).Value = "500,000 円 ~ 1,000,000 円 / 固定"
$ws.Cells.Item(6, 5).Value = "期限情報なし"
$url6 = "https://www.lancers.jp/work/detail/5460724"
$ws.Cells.Item(6, 6).Value = $url6
$ws.Cells.Item(6, 7).Value = 125
$ws.Cells.Item(6, 8).Value = "◆開発,システム開発"

# Row 7 - new listing (was previously at row 10).
$ws.Cells.Item(7, 2).Value = "【急募】宿泊業向けSaaSの予約者取得システム開発"
$ws.Cells.Item(7, 3).Value = "システム開発"
$ws.Cells.Item(7, 4).Value = "100,000 円 ~ 200,000 円 / 固定"
$ws.Cells.Item(7, 5).Value = "期限情報なし"
$url7 = "https://www.lancers.jp/work/detail/5460405"
$ws.Cells.Item(7, 6).Value = $url7
$ws.Cells.Item(7, 7).Value = 118
$ws.Cells.Item(7, 8).Value = "◆開発,システム開発"

# Row 8 - brand new listing.
$ws.Cells.Item(8, 2).Value = "【急募】顧客管理システムの開発をお手伝いください!"
$ws.Cells.Item(8, 3).Value = "システム開発"
$ws.Cells.Item(8, 4).Value = "500,000 円 ~ 1,000,000 円 / 固定"
$ws.Cells.Item(8, 5).Value = "期限情報なし"
$url8 = "https://www.lancers.jp/work/detail/5460928"
$ws.Cells.Item(8, 6).Value = $url8
$ws.Cells.Item(8, 7).Value = 115
$ws.Cells.Item(8, 8).Value = "◆開発 ◇管理"

# Row 9 - brand new listing.
$ws.Cells.Item(9, 2).Value = "【急募】PHPシステムエラーの早期改善を求む!"
$ws.Cells.Item(9, 3).Value = "システム開発"
$ws.Cells.Item(9, 4).Value = "10,000 円 ~ 20,000 円 / 固定"
$ws.Cells.Item(9, 5).Value = "期限情報なし"
$url9 = "https://www.lancers.jp/work/detail/5460787"
$ws.Cells.Item(9, 6).Value = $url9
$ws.Cells.Item(9, 7).Value = 40
$ws.Cells.Item(9, 8).Value = "○PHP"

# Row 10 - brand new listing.
$ws.Cells.Item(10, 2).Value = "wordpressレンダリングを妨げるリソースの除外"
$ws.Cells.Item(10, 3).Value = "システム開発"
$ws.Cells.Item(10, 4).Value = "200,000 円 ~ 300,000 円 / 固定"
$ws.Cells.Item(10, 5).Value = "期限情報なし"
$url10 = "https://www.lancers.jp/work/detail/5016989"
$ws.Cells.Item(10, 6).Value = $url10
$ws.Cells.Item(10, 7).Value = 33
$ws.Cells.Item(10, 8).Value = "○WordPress"

# Row 11 - new listing (was previously at row 14).
$ws.Cells.Item(11, 2).Value = "【急募】WEBサイト研修講師を探しています!"
$ws.Cells.Item(11, 3).Value = "システム開発"
$ws.Cells.Item(11, 4).Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Cells.Item(11, 5).Value = "期限情報なし"
$url11 = "https://www.lancers.jp/work/detail/5460484"
$ws.Cells.Item(11, 6).Value = $url11
$ws.Cells.Item(11, 7).Value = 33
$ws.Cells.Item(11, 8).Value = "◇サイト"

# Row 12 - new listing (was previously at row 19); this one has no
# skill-summary column, so H12 must end up empty.
$ws.Cells.Item(12, 2).Value = "限定公開 限定公開の仕事"
$ws.Cells.Item(12, 3).Value = "システム開発"
$ws.Cells.Item(12, 4).Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Cells.Item(12, 5).Value = "期限情報なし"
$url12 = "https://www.lancers.jp/work/detail/5450323"
$ws.Cells.Item(12, 6).Value = $url12
$ws.Cells.Item(12, 7).Value = 13
$ws.Cells.Item(12, 8).ClearContents()

# ---------------------------------------------------------------------
# 5) Rebuild the F2:F12 hyperlinks so each one targets the URL that is
#    now displayed in that cell.
# ---------------------------------------------------------------------
$ws.Hyperlinks.Add($ws.Cells.Item(2, 6), $url2)
$ws.Hyperlinks.Add($ws.Cells.Item(3, 6), $url3)
$ws.Hyperlinks.Add($ws.Cells.Item(4, 6), $url4)
$ws.Hyperlinks.Add($ws.Cells.Item(5, 6), $url5)
$ws.Hyperlinks.Add($ws.Cells.Item(6, 6), $url6)
$ws.Hyperlinks.Add($ws.Cells.Item(7, 6), $url7)
$ws.Hyperlinks.Add($ws.Cells.Item(8, 6), $url8)
$ws.Hyperlinks.Add($ws.Cells.Item(9, 6), $url9)
$ws.Hyperlinks.Add($ws.Cells.Item(10, 6), $url10)
$ws.Hyperlinks.Add($ws.Cells.Item(11, 6), $url11)
$ws.Hyperlinks.Add($ws.Cells.Item(12, 6), $url12)

# ---------------------------------------------------------------------
# 6) Column D narrows from 32 to 30 characters. ColumnWidth round-trips
#    through a pixel-snapped unit, so 29.1667 is the input that lands on
#    a stored width of exactly 30.
# ---------------------------------------------------------------------
$ws.Columns.Item(4).ColumnWidth = 29.1667
